# Clean up the data processing output: the profiling metadata previously
# stored the literal column "position" and a best-guess "type" per
# attribute. Those are replaced so that position is no longer tracked
# (-1 placeholder) and the type/derived statistics reflect that the
# value is now produced generically instead of being hard-coded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: person_id_nbr ---
$ws.Range("B2").Value = -1
$ws.Range("C2").Value = "integer"

# --- Row 3: surname ---
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "date"

# --- Row 4: given_name ---
$ws.Range("B4").Value = -1
$ws.Range("C4").Value = "date"
$ws.Range("U4").Value = 0

# --- Row 5: year ---
$ws.Range("B5").Value = -1
$ws.Range("C5").Value = "date"
$ws.Range("E5").Value = 83260
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -1
$ws.Range("L5").Value = -1

# --- Row 6: place ---
$ws.Range("B6").Value = -1
$ws.Range("C6").Value = "date"

# --- Row 7: title ---
$ws.Range("B7").Value = -1
$ws.Range("C7").Value = "date"
$ws.Range("U7").Value = 0

# --- Row 8: reference ---
$ws.Range("B8").Value = -1
$ws.Range("C8").Value = "date"
$ws.Range("U8").Value = 0

# --- Row 9: volume ---
$ws.Range("B9").Value = -1
$ws.Range("C9").Value = "date"
$ws.Range("U9").Value = 0

# --- Row 10: bundle ---
$ws.Range("B10").Value = -1
$ws.Range("C10").Value = "date"

# --- Row 11: petition ---
$ws.Range("B11").Value = -1
$ws.Range("C11").Value = "date"
$ws.Range("E11").Value = 83260
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = -1
$ws.Range("L11").Value = -1

# --- Row 12: microfilm ---
$ws.Range("B12").Value = -1
$ws.Range("C12").Value = "date"
$ws.Range("U12").Value = 0

# --- Row 13: page ---
$ws.Range("B13").Value = -1
$ws.Range("C13").Value = "date"
